# Gallery_PanelNode_Repeaters_ProPanels.xlsx
# "Added and Updated Czech Test Data"
#
# 1) Duplicate the "Belgium" sheet (same layout/styles/merges) right after it,
#    rename the copy to "Czech".
# 2) Update the market name / ticket reference cells on the new sheet.
# 3) Remove the two repeater rows ("ZXF" / "ZXFEV") that Czech does not test.
# 4) Resize columns B:D on the new sheet to the Czech-specific widths.
# 5) Fix up selections: Belgium's selection becomes the full used range (no
#    longer the active tab), Czech becomes the active tab with D19 selected.

$wb = $excel.ActiveWorkbook

$belgium = $wb.Worksheets.Item("Belgium")

# Belgium's on-screen selection after the edit spans its whole used range.
$belgium.Range("A1:D23").Select() | Out-Null

# Duplicate Belgium right after itself -> becomes the 3rd sheet.
$belgium.Copy($null, $belgium)
$czech = $wb.Worksheets.Item(3)
$czech.Name = "Czech"

# Market name (row 2) and Jira reference (row 4) for the Czech sheet.
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1734"

# Czech doesn't cover the "ZXF" / "ZXFEV" repeaters - drop those two rows.
$czech.Rows("20:21").Delete()

# Column widths specific to the Czech sheet.
$czech.Columns("B").ColumnWidth = 22.8
$czech.Columns("C").ColumnWidth = 16.2
$czech.Columns("D").ColumnWidth = 29

# Final selection / active sheet.
$czech.Range("D19").Select() | Out-Null
$czech.Activate()
